$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 708/709 (Primera + Segunda for a new weekly
# price report date), pushing the existing rows 708:743 down to 710:745.
$ws.Rows.Item(708).Resize(2).EntireRow.Insert()

# Row 708: new "Primera" quality record
$ws.Range("A708").Value2 = 8
$ws.Range("B708").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C708").Value2 = "Coquimbo"
$ws.Range("D708").Value2 = 44753
$ws.Range("E708").Value2 = 4
$ws.Range("F708").Value2 = 100112023
$ws.Range("G708").Value2 = "Brócoli"
$ws.Range("H708").Value2 = "Sin especificar"
$ws.Range("I708").Value2 = "Primera"
$ws.Range("J708").Value2 = 2600
$ws.Range("K708").Value2 = 800
$ws.Range("L708").Value2 = 900
$ws.Range("M708").Value2 = 850
$ws.Range("N708").Value2 = "$/unidad"
$ws.Range("O708").Value2 = "Provincia del Elquí"
$ws.Range("P708").Value2 = 850
$ws.Range("Q708").Value2 = 1
$ws.Range("R708").Value2 = "Hortaliza"

# Row 709: new "Segunda" quality record (same date/market as row 708)
$ws.Range("A709").Value2 = 8
$ws.Range("B709").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C709").Value2 = "Coquimbo"
$ws.Range("D709").Value2 = 44753
$ws.Range("E709").Value2 = 4
$ws.Range("F709").Value2 = 100112023
$ws.Range("G709").Value2 = "Brócoli"
$ws.Range("H709").Value2 = "Sin especificar"
$ws.Range("I709").Value2 = "Segunda"
$ws.Range("J709").Value2 = 1360
$ws.Range("K709").Value2 = 700
$ws.Range("L709").Value2 = 750
$ws.Range("M709").Value2 = 725
$ws.Range("N709").Value2 = "$/unidad"
$ws.Range("O709").Value2 = "Provincia del Elquí"
$ws.Range("P709").Value2 = 725
$ws.Range("Q709").Value2 = 1
$ws.Range("R709").Value2 = "Hortaliza"
